$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove workbook protection (workbookProtection element)
$wb.Unprotect()

# Update row 2 values: A2 becomes text "7034", B2 becomes text "98745632"
$ws.Range("A2").Value = "'7034"
$ws.Range("B2").Value = "'98745632"

# Remove row 3 entirely (shrinks used range to A1:B2)
$ws.Rows.Item(3).Delete()

# Reproduce the saved selection/active-window state from the diff
[void]$ws.Range("A2:XFD4").Select()
